{"js": "// The campaign blurb \"Dates \u00e0 utiliser pour la Campagne Constellation du\n// Taureau 2022: 16-25 janvier\" is repeated (identically) several times\n// throughout the document body. Every occurrence needs to be updated so\n// that \"2022\" moves right after \"Campagne\" instead of trailing the\n// constellation name:\n//   \"Campagne Constellation du Taureau 2022\" -> \"Campagne 2022 Constellation du Taureau\"\nconst oldText = \"Dates \u00e0 utiliser pour la Campagne Constellation du Taureau 2022: 16-25 janvier\";\nconst newText = \"Dates \u00e0 utiliser pour la Campagne 2022 Constellation du Taureau: 16-25 janvier\";\n\nconst results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The campaign blurb \"Dates \u00e0 utiliser pour la Campagne Constellation du\n# Taureau 2022: 16-25 janvier\" occurs several times (identically) throughout\n# the document body. Every occurrence must be updated so that \"2022\" moves\n# right after \"Campagne\" instead of trailing the constellation name:\n#   \"Campagne Constellation du Taureau 2022\" -> \"Campagne 2022 Constellation du Taureau\"\n#\n# A single Find/Replace (Replace:=wdReplaceAll) over the whole document\n# content range takes care of every occurrence in one pass.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"Dates \u00e0 utiliser pour la Campagne Constellation du Taureau 2022: 16-25 janvier\"\n$find.Replacement.Text = \"Dates \u00e0 utiliser pour la Campagne 2022 Constellation du Taureau: 16-25 janvier\"\n\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n# wdReplaceAll = 2 replaces every match found in the range.\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
